# Fruta / hortaliza, semanal
# Insert a new daily-price record row at row 288 (pushing the existing
# rows 288-380 down to 289-381) and populate it with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 288; Excel shifts rows 288..380 down to 289..381
# and carries the surrounding number formatting (e.g. the date style on D).
$ws.Rows.Item(288).Insert()

$ws.Range("A288").Value = 10
$ws.Range("B288").Value = "Vega Modelo de Temuco"
$ws.Range("C288").Value = "La Araucanía"
$ws.Range("D288").Value = 44524
$ws.Range("E288").Value = 9
$ws.Range("F288").Value = "Fruta"
$ws.Range("G288").Value = 100108
$ws.Range("H288").Value = "Tropicales y subtropicales"
$ws.Range("I288").Value = 100108006
$ws.Range("J288").Value = "Plátano"
$ws.Range("K288").Value = "Sin especificar"
$ws.Range("L288").Value = "Pintón"
$ws.Range("M288").Value = 730
$ws.Range("N288").Value = 23000
$ws.Range("O288").Value = 24000
$ws.Range("P288").Value = 23479
$ws.Range("Q288").Value = "$/caja 20 kilos"
$ws.Range("R288").Value = "Ecuador"
$ws.Range("S288").Value = 1174
$ws.Range("T288").Value = 20
